$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.86"
$ws.Range("E2").Value = "'-0.74%"
$ws.Range("G2").Value = "'17"

$ws.Range("D3").Value = "'27.01"
$ws.Range("E3").Value = "'2.08%"
$ws.Range("G3").Value = "'17"

$ws.Range("D4").Value = "'5.073"
$ws.Range("E4").Value = "'-0.15%"
$ws.Range("G4").Value = "'17"

$ws.Range("D5").Value = "'0.05688"
$ws.Range("G5").Value = "'17"

$ws.Range("D6").Value = "'6.477"
$ws.Range("E6").Value = "'-0.18%"
$ws.Range("G6").Value = "'17"

$ws.Range("D7").Value = "'0.8204"
$ws.Range("E7").Value = "'0.79%"
$ws.Range("G7").Value = "'17"

$ws.Range("D8").Value = "'0.8393"
$ws.Range("E8").Value = "'-0.63%"
$ws.Range("G8").Value = "'17"

$ws.Range("D9").Value = "'0.1326"
$ws.Range("E9").Value = "'-1.43%"
$ws.Range("G9").Value = "'17"

$ws.Range("D10").Value = "'0.06890"
$ws.Range("E10").Value = "'-1.18%"
$ws.Range("G10").Value = "'17"

$ws.Range("D11").Value = "'0.02855"
$ws.Range("E11").Value = "'-0.62%"
$ws.Range("G11").Value = "'17"

$ws.Range("E12").Value = "'-0.15%"
$ws.Range("G12").Value = "'17"

$ws.Range("D13").Value = "'0.001508"
$ws.Range("E13").Value = "'-1.29%"
$ws.Range("G13").Value = "'17"

$ws.Range("D14").Value = "'0.04089"
$ws.Range("E14").Value = "'-12.44%"
$ws.Range("G14").Value = "'17"

$ws.Range("D15").Value = "'0.0006019"
$ws.Range("E15").Value = "'0.88%"
$ws.Range("G15").Value = "'17"

$ws.Range("D16").Value = "'0.006073"
$ws.Range("E16").Value = "'-2.09%"
$ws.Range("G16").Value = "'17"

$ws.Range("E17").Value = "'-2.30%"
$ws.Range("G17").Value = "'17"

$ws.Range("D18").Value = "'3.001"
$ws.Range("E18").Value = "'-0.29%"
$ws.Range("G18").Value = "'17"

$ws.Range("D19").Value = "'2.226"
$ws.Range("E19").Value = "'5.08%"
$ws.Range("G19").Value = "'17"

$ws.Range("E20").Value = "'-0.22%"
$ws.Range("G20").Value = "'17"

$ws.Range("D21").Value = "'0.03170"
$ws.Range("E21").Value = "'-0.36%"
$ws.Range("G21").Value = "'17"

$ws.Range("E22").Value = "'-1.79%"
$ws.Range("G22").Value = "'17"

$ws.Range("D23").Value = "'3.586"
$ws.Range("E23").Value = "'-4.74%"
$ws.Range("G23").Value = "'17"

$ws.Range("E24").Value = "'1.80%"
$ws.Range("G24").Value = "'17"

$ws.Range("D25").Value = "'0.001219"
$ws.Range("E25").Value = "'-2.34%"
$ws.Range("G25").Value = "'17"

$ws.Range("D26").Value = "'0.003951"
$ws.Range("E26").Value = "'-14.15%"
$ws.Range("G26").Value = "'17"

$ws.Range("G27").Value = "'17"

$ws.Range("D28").Value = "'0.0001937"
$ws.Range("E28").Value = "'0.02%"
$ws.Range("G28").Value = "'17"

$ws.Range("G29").Value = "'17"

$ws.Range("G30").Value = "'17"

$ws.Range("G31").Value = "'17"

$ws.Range("G32").Value = "'17"

$ws.Range("G33").Value = "'17"

$ws.Range("G34").Value = "'17"

$ws.Range("G35").Value = "'17"

$ws.Range("G36").Value = "'17"

$ws.Range("G37").Value = "'17"

$ws.Range("G38").Value = "'17"

$ws.Range("G39").Value = "'17"

$ws.Range("D40").Value = "'0.03688"
$ws.Range("E40").Value = "'0.14%"
$ws.Range("G40").Value = "'17"

$ws.Range("D41").Value = "'0.005877"
$ws.Range("E41").Value = "'-5.06%"
$ws.Range("G41").Value = "'17"

$ws.Range("D42").Value = "'0.1055"
$ws.Range("E42").Value = "'-0.42%"
$ws.Range("G42").Value = "'17"

$ws.Range("D43").Value = "'0.002335"
$ws.Range("E43").Value = "'-6.62%"
$ws.Range("G43").Value = "'17"

$ws.Range("E44").Value = "'5.17%"
$ws.Range("G44").Value = "'17"

$ws.Range("D45").Value = "'0.00005216"
$ws.Range("E45").Value = "'-1.42%"
$ws.Range("G45").Value = "'17"

$ws.Range("E46").Value = "'-0.03%"
$ws.Range("G46").Value = "'17"

$ws.Range("E47").Value = "'-32.27%"
$ws.Range("G47").Value = "'17"

$ws.Range("E48").Value = "'2.66%"
$ws.Range("G48").Value = "'17"

$ws.Range("E49").Value = "'-0.03%"
$ws.Range("G49").Value = "'17"

$ws.Range("E50").Value = "'-0.03%"
$ws.Range("G50").Value = "'17"

$ws.Range("G51").Value = "'17"
